$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("H3").Value = 0.8878364905284147
$ws.Range("I3").Value = 0.007966602577796767
$ws.Range("K3").Value = 334.6451612903226

$ws.Range("Q3").Value = 9
$ws.Range("R3").Value = 27
$ws.Range("S3").Value = 104
$ws.Range("T3").Value = 225
$ws.Range("U3").Value = 632
$ws.Range("V3").Value = 38043
$ws.Range("W3").Value = 38025
$ws.Range("X3").Value = 37948
$ws.Range("Y3").Value = 37827
$ws.Range("Z3").Value = 37420

$ws.Range("AF3").Value = 0.999763
$ws.Range("AG3").Value = 0.99929
$ws.Range("AH3").Value = 0.997267
$ws.Range("AI3").Value = 0.9940870000000001
$ws.Range("AJ3").Value = 0.983391
